# Applies the "Close to a final draft" edit:
#  1. Refreshes the cached datetimeFigureOut field text (10/12/2018 -> 19/12/2018)
#     on the slide master and every slide layout.
#  2. Repositions / resizes the three pictures and three "(a)/(b)/(c)" caption
#     textboxes on slide 1, and crops the left edge of Picture 2.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder field text
# ---------------------------------------------------------------------------
function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "10/12/2018") {
                    $shp.TextFrame.TextRange.Text = "19/12/2018"
                }
            }
        }
    }
}

# Slide master
Update-DateField $p.SlideMaster.Shapes

# Every slide layout
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateField $layouts.Item($j).Shapes
}

# ---------------------------------------------------------------------------
# 2. Slide 1 shape geometry
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# Picture 4 (big right-hand screenshot)
$picture4 = $s.Shapes.Item("Picture 4")
$picture4.Left   = 243.09047744094488
$picture4.Top    = -1.2203942007874016
$picture4.Width  = 313.9943392086614
$picture4.Height = 458.2979587559055

# TextBox 1 -> "(a)" label
$textBoxA = $s.Shapes.Item("TextBox 1")
$textBoxA.Left = 178.66551181102363
$textBoxA.Top  = 167.64968503937007

# TextBox 6 -> "(b)" label
$textBoxB = $s.Shapes.Item("TextBox 6")
$textBoxB.Left = 106.80078740157481
$textBoxB.Top  = 456.7830051259843

# TextBox 7 -> "(c)" label
$textBoxC = $s.Shapes.Item("TextBox 7")
$textBoxC.Left = 380.71031496062994
$textBoxC.Top  = 456.7830051259843

# Picture 2 (top-left picture) - cropped on the left edge and moved/resized
$picture2 = $s.Shapes.Item("Picture 2")
$picture2.LockAspectRatio = 0
$picture2.PictureFormat.CropLeft = 19.8324
$picture2.Left   = 9.76251968503937
$picture2.Top    = 0.0
$picture2.Width  = 173.2224409448819
$picture2.Height = 203.58409898818897

# Picture 5 (bottom-left picture) - moved only
$picture5 = $s.Shapes.Item("Picture 5")
$picture5.Left = 9.76251968503937
$picture5.Top  = 216.6215748031496

Write-Host "edit applied"
